$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nicolas Blavet (row 6) ---
# D6: "Site perso" hyperlink
$ws.Hyperlinks.Add($ws.Range("D6"), "http://nikola7654.wixsite.com/monsite")
$ws.Range("D6").Style = "Lien hypertexte"

# E6 (CV column): highlighted, empty cell (greenish text/fill, no hyperlink)
$ws.Range("E6").Interior.Color = 13499135
$ws.Range("E6").Font.ThemeColor = 9
$ws.Range("E6").Font.TintAndShade = 0.39997558519241921

# G6: "Linkdin" hyperlink
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.linkedin.com/in/nicolas-blavet-666285136/")
$ws.Range("G6").Style = "Lien hypertexte"

# --- Melinda Khammar (row 4) ---
# C4 (Photo column): highlighted, empty cell (just the fill)
$ws.Range("C4").Interior.Color = 13499135

# D4: "Site perso" hyperlink, keep the bordered/filled look of the row
$ws.Hyperlinks.Add($ws.Range("D4"), "http://www.mekabull.fr/")
$ws.Range("D4").Style = "Lien hypertexte"
$ws.Range("D4").Interior.Color = 15921906

# F4: twitter handle as plain text (no hyperlink relationship)
$ws.Range("F4").Value = "https://twitter.com/mekabulle"

# G4: Linkedin profile as plain text (no hyperlink relationship)
$ws.Range("G4").Value = "https://fr.linkedin.com/in/melinda-khammar-868885133"

$ws.Range("D25").Select()
